$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# New row of regression results for the "Baseline 2010-18" run
$ws.Range("A7").Value = "CW3M"
$ws.Range("B7").Value = "Baseline_2010-18_c45 9/19/20"
$ws.Range("C7").Value = "2010-18"
$ws.Range("D7").Value = 1138.6194117777777
$ws.Range("E7").Value = 1901.5157334444443
$ws.Range("F7").Value = 1.0119255555555557
$ws.Range("G7").Value = 327.78053433333326
$ws.Range("H7").Value = 9.775355222222224
$ws.Range("I7").Value = 8.145128999999999
$ws.Range("J7").Value = 769.26639155555551
$ws.Range("K7").Value = 83.47062044444445
$ws.Range("L7").Value = 1374.8233372222221
$ws.Range("M7").Value = 1142.9502087777778
$ws.Range("N7").Value = 4918.1879612222219
$ws.Range("O7").Value = 27227.338324777778
$ws.Range("P7").Value = -0.0472741111111111
$ws.Range("Q7").Value = -0.000038888888888888877
$ws.Range("R7").Value = "2010-18"

# Match number formats used by the existing data rows
$ws.Range("D7:M7").NumberFormat = "0.00"
$ws.Range("N7:O7").NumberFormat = "0"
$ws.Range("P7:Q7").NumberFormat = "0.000000"

$ws.Range("R8").Select()
